$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title row (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 20:05"

# Update country rows: name (col A) and stats (cols B-H)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1696475
$ws.Range("C4").Value = 10039
$ws.Range("D4").Value = 456595
$ws.Range("E4").Value = 1140324
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 256
$ws.Range("H4").Value = 99556

$ws.Range("A13").Value = "India"
$ws.Range("B13").Value = 144868
$ws.Range("C13").Value = 6332
$ws.Range("D13").Value = 60698
$ws.Range("E13").Value = 79998
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 148
$ws.Range("H13").Value = 4172

$ws.Range("A32").Value = "Emiratos Arabes Unidos"
$ws.Range("B32").Value = 30307
$ws.Range("C32").Value = 822
$ws.Range("D32").Value = 15657
$ws.Range("E32").Value = 14402
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 248

$ws.Range("A33").Value = "Irlanda"
$ws.Range("B33").Value = 24698
$ws.Range("C33").Value = 59
$ws.Range("D33").Value = 21060
$ws.Range("E33").Value = 2032
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 1606

$ws.Range("A118").Value = "Paraguay"
$ws.Range("B118").Value = 865
$ws.Range("C118").Value = 3
$ws.Range("D118").Value = 344
$ws.Range("E118").Value = 510
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 11

$ws.Range("A131").Value = "Republica de Africa Central"
$ws.Range("B131").Value = 652
$ws.Range("C131").Value = 48
$ws.Range("D131").Value = 22
$ws.Range("E131").Value = 629
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 1

$ws.Range("A132").Value = "Malta"
$ws.Range("B132").Value = 611
$ws.Range("C132").Value = 1
$ws.Range("D132").Value = 485
$ws.Range("E132").Value = 120
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 6

$ws.Range("A151").Value = "Suazilandia"
$ws.Range("B151").Value = 256
$ws.Range("C151").Value = 6
$ws.Range("D151").Value = 158
$ws.Range("E151").Value = 96
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 2

$ws.Range("A198").Value = "Curazao"
$ws.Range("B198").Value = 18
$ws.Range("C198").Value = 1
$ws.Range("D198").Value = 14
$ws.Range("E198").Value = 3
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Santa Lucia"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 18
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

$ws.Range("A200").Value = "Nueva Caledonia"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 16
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 2

$ws.Range("A202").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B202").Value = 17
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 0
$ws.Range("E202").Value = 17
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 10
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
